$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "-"
$ws.Range("C4").Value = "-"
$ws.Range("E4").Value = "-"
$ws.Range("B6").Value = "['MCT-2A-Sistemas digitais', -, -]"
$ws.Range("E6").Value = "-"
$ws.Range("B7").Value = "['MCT-2A-Sistemas digitais', -, -]"
$ws.Range("C7").Value = "[-, -, 'MCT-3A-Lab. de eletroeletrônica']"
$ws.Range("C8").Value = "[-, -, 'MCT-3A-Lab. de eletroeletrônica']"
